$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Delete columns J:K (they held the now-removed "tendency"/"OBSERVED" helper column)
$ws.Range("J1:K1").EntireColumn.Delete()

# Clear the stray space left in D8
$ws.Range("D8").ClearContents()

# Delete rows 9:12 (blank filler rows removed, shifting the footer rows up)
$ws.Range("A9:A12").EntireRow.Delete()

# Refresh the cached sort-state range so it tracks the new (shrunk) data block
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A5:A9"))
$ws.Sort.SetRange($ws.Range("A5:A9"))
$ws.Sort.Apply()

# Update the active selection to match the post-edit state
$ws.Range("K9").Select() | Out-Null
